$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. New styles: "Procedure Intro" (paragraph) / "Procedure Intro Char"
#    (linked character style) and the "Hyperlink" character style.
# ---------------------------------------------------------------------------

$pi = $d.Styles.Add("ProcedureIntro", 1)            # wdStyleTypeParagraph
$pi.NameLocal = "Procedure Intro"
$pi.Font.Name = "Arial"
$pi.Font.NameFarEast = "Times New Roman"
$pi.Font.NameBi = "Arial"
$pi.Font.Bold = $true
$pi.Font.BoldBi = $true
$pi.Font.Size = 10
$pi.Font.SizeBi = 12
$pi.NextParagraphStyle = "Normal"
$pi.ParagraphFormat.SpaceBefore = 6
$pi.ParagraphFormat.SpaceAfter = 6
$pi.ParagraphFormat.LineSpacingRule = 0            # wdLineSpaceSingle

$pic = $d.Styles.Add("ProcedureIntroChar", 2)       # wdStyleTypeCharacter
$pic.NameLocal = "Procedure Intro Char"
$pic.BaseStyle = "DefaultParagraphFont"
$pic.Font.Name = "Arial"
$pic.Font.NameFarEast = "Times New Roman"
$pic.Font.NameBi = "Arial"
$pic.Font.Bold = $true
$pic.Font.BoldBi = $true
$pic.Font.Size = 10
$pic.Font.SizeBi = 12

$pi.LinkStyle = "ProcedureIntroChar"
$pic.LinkStyle = "ProcedureIntro"

$hl = $d.Styles.Add("Hyperlink", 2)                 # wdStyleTypeCharacter
$hl.BaseStyle = "DefaultParagraphFont"
$hl.Priority = 99
$hl.UnhideWhenUsed = $true
$hl.Font.Color = 16711680                           # OLE BGR -> w:color 0000FF
$hl.Font.Underline = 1                               # wdUnderlineSingle

# ---------------------------------------------------------------------------
# 2. Second paragraph (hyperlink) and third (empty trailing) paragraph are
#    built first, while paragraph 1 is still in its pristine state, so that
#    neither inherits the ProcedureIntro style/formatting applied below.
# ---------------------------------------------------------------------------

$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$p2 = $d.Paragraphs($d.Paragraphs.Count)
$p2Range = $p2.Range
$p2Range.Collapse(0)
$p2Range.InsertAfter(" ")

$hlInsertRng = $p2.Range
$hlInsertRng.Collapse(1)   # wdCollapseStart
$d.Hyperlinks.Add($hlInsertRng, "http://www.vmware.com/go/vFabric-ref-arch", $null, $null, "http://www.vmware.co/go/vFabric-ref-Arch") | Out-Null

$finalRng = $d.Content
$finalRng.Collapse(0)
$finalRng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# ---------------------------------------------------------------------------
# 3. First paragraph: replace "vFRA" (+ bookmark) with the new intro text
#    and apply the ProcedureIntro paragraph style.
# ---------------------------------------------------------------------------

$p1 = $d.Paragraphs(1).Range
$introXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Materials for Topic 3</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t>&#8220;</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t>Capacity Planning</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t>&#8221;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> and all vFabric Reference Architecture topics are located at: </w:t>
  </w:r>
</w:p>
'@
$p1.InsertXML($introXml)
$d.Paragraphs(1).Style = "ProcedureIntro"

Write-Host "done"
